$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Create the two new lookup sheets, right after "comportamiento".
#    causa_probable is created first (so it receives the lower internal
#    sheetId), then tipo_registro is created (also anchored right after
#    comportamiento, so it ends up positioned BEFORE causa_probable).
# ---------------------------------------------------------------------------
$anchor = $wb.Worksheets.Item("comportamiento")
$wsCausaNew = $wb.Worksheets.Add($null, $anchor)
$wsCausaNew.Name = "causa_probable"
$wsTipoNew = $wb.Worksheets.Add($null, $anchor)
$wsTipoNew.Name = "tipo_registro"

# Re-fetch stable references by name (Add() handles can alias after a
# second Add() changes relative sheet positions).
$wsCausa = $wb.Worksheets.Item("causa_probable")
$wsTipo = $wb.Worksheets.Item("tipo_registro")

# ---------------------------------------------------------------------------
# 2) Populate causa_probable (A1:A8) -- written first so its new shared
#    strings are allocated before tipo_registro's.
# ---------------------------------------------------------------------------
$causaRows = @(
    "causa_probable",
    "Atropello",
    "Natural",
    "Por ataque",
    "Ahogado",
    "Envenenado",
    "Colision",
    "Otra causa"
)
for ($i = 0; $i -lt $causaRows.Length; $i++) {
    $wsCausa.Cells.Item($i + 1, 1).Value2 = $causaRows[$i]
}

# ---------------------------------------------------------------------------
# 3) Populate tipo_registro (A1:B29)
# ---------------------------------------------------------------------------
$tipoRows = @(
    @("tipo_registro","componente_biologico"),
    @("Auditivo","Ornitofauna"),
    @("Avistamiento","Herpetofauna"),
    @("Avistamiento","Mastofauna"),
    @("Avistamiento","Ornitofauna"),
    @("Captura temporal","Mastofauna"),
    @("Captura temporal","Herpetofauna"),
    @("Captura temporal","Ornitofauna"),
    @("Carcasa/Restos/Osamenta","Mastofauna"),
    @("Carcasa/Restos/Osamenta","Herpetofauna"),
    @("Carcasa/Restos/Osamenta","Ornitofauna"),
    @("Dormidero","Ornitofauna"),
    @("Dormidero","Mastofauna"),
    @("Egagropila","Ornitofauna"),
    @("Estercolero","Mastofauna"),
    @("Excavacion","Mastofauna"),
    @("Feca","Mastofauna"),
    @("Feca","Herpetofauna"),
    @("Feca","Ornitofauna"),
    @("Galerias","Mastofauna"),
    @("Huella","Mastofauna"),
    @("Huella","Herpetofauna"),
    @("Huella","Ornitofauna"),
    @("Nido","Ornitofauna"),
    @("Pelos","Mastofauna"),
    @("Piel","Mastofauna"),
    @("Piel","Herpetofauna"),
    @("Plumas","Ornitofauna"),
    @("Revolcadero","Mastofauna")
)
for ($i = 0; $i -lt $tipoRows.Length; $i++) {
    $r = $i + 1
    $wsTipo.Cells.Item($r, 1).Value2 = $tipoRows[$i][0]
    $wsTipo.Cells.Item($r, 2).Value2 = $tipoRows[$i][1]
}

# ---------------------------------------------------------------------------
# 4) Turn both ranges into Tables (ListObjects), matching names used by the
#    source workbook ("Tabla4" for causa_probable, "Tabla10" for tipo_registro).
# ---------------------------------------------------------------------------
$loCausa = $wsCausa.ListObjects.Add(1, $wsCausa.Range("A1:A8"), $null, 1)
$loCausa.Name = "Tabla4"

$loTipo = $wsTipo.ListObjects.Add(1, $wsTipo.Range("A1:B29"), $null, 1)
$loTipo.Name = "Tabla10"

# ---------------------------------------------------------------------------
# 5) Tab colors for the new sheets (match the green used by the other
#    lookup sheets).
# ---------------------------------------------------------------------------
$wsCausa.Tab.Color = 5296274
$wsTipo.Tab.Color = 5296274

# ---------------------------------------------------------------------------
# 6) Column widths -- autofit based on the content just entered.
# ---------------------------------------------------------------------------
$wsCausa.Columns.Item(1).AutoFit() | Out-Null
$wsTipo.Columns.Item(1).AutoFit() | Out-Null
$wsTipo.Columns.Item(2).AutoFit() | Out-Null

# ---------------------------------------------------------------------------
# 7) Update "comportamiento": tab color, column A width, becomes the
#    active/selected sheet with a new selection.
# ---------------------------------------------------------------------------
$wsComp = $wb.Worksheets.Item("comportamiento")
$wsComp.Tab.Color = 5296274
$wsComp.Columns.Item(1).AutoFit() | Out-Null

# ---------------------------------------------------------------------------
# 8) Selections: mot_interv (previously active) moves its selection and
#    loses the "active" flag; comportamiento becomes active with B17
#    selected.
# ---------------------------------------------------------------------------
$wsMot = $wb.Worksheets.Item("mot_interv")
$wsMot.Activate()
$wsMot.Range("B27").Select() | Out-Null

$wsComp.Activate()
$wsComp.Range("B17").Select() | Out-Null
